$d = $word.ActiveDocument

# Locate the "Test" run using Find (as recommended), then re-derive a
# plain Range from its Start/End so the subsequent .Text assignment
# splits the paragraph cleanly and lets the trailing _GoBack bookmark
# settle into its own (new, final) paragraph - mirroring what Word does
# when a user places the caret right after "Test" and presses Enter a
# few times before typing the new sentence.
$r = $d.Content
$found = $r.Find.Execute("Test", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $start = $r.Start
    $end = $r.End
    $target = $d.Range($start, $end)
    $target.Text = "Test`r`rJ’ajoute une ligne pour m’exercer aux commandes git`r`r"
}
